$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.088048696517944
$ws.Range("B1").Value = 1.168573975563049
$ws.Range("C1").Value = 1.120832324028015
$ws.Range("D1").Value = 1.33975076675415
$ws.Range("E1").Value = 1.253533124923706
